# Apply the "Add files via upload" edit to 360MasterData.xlsx
# Summary of the change:
#  - A new scene row (row 24) is appended to the Scenes table, re-using the
#    data that used to live in F2 (PanoramaURL) together with a Voice-Over
#    hyperlink, a Logo URL hyperlink and a Background-music-URL hyperlink.
#  - F2 (PanoramaURL for the first scene) is updated to point at a brand new
#    HDRI URL and its "Background music Active" flag (K2) is switched on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1. Scenes Sheet")

# ---------------------------------------------------------------------
# 1. Append a new row (24) at the bottom of Table2, copying row 23's
#    data/format, then tweak the handful of cells that differ.
# ---------------------------------------------------------------------
$srcRow = $ws.Range("A23:R23")
$dstRow = $ws.Range("A24:R24")
$srcRow.Copy()
$dstRow.PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# Column F (PanoramaURL) on the new row takes the URL that used to be in F2
$ws.Range("F24").Value = "https://i.imgur.com/58TUihy.jpeg"
$ws.Hyperlinks.Add($ws.Range("F24"), "https://i.imgur.com/58TUihy.jpeg") | Out-Null

# Column Q (SceneID Active) turns on for the new row
$ws.Range("Q24").Value = $true

# Re-assert the Logo URL / Background music URL hyperlinks (kept the same
# targets as the row above, but Excel records a fresh hyperlink entry for
# them because they were pasted into a brand new row)
$ws.Hyperlinks.Add($ws.Range("G24"), "https://i.imgur.com/6DLBULh.jpeg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H24"), "https://raw.githubusercontent.com/Rmoosa2014/vr-tour/main/BGM.mp3") | Out-Null

# ---------------------------------------------------------------------
# 2. Update F2 (PanoramaURL for the very first scene) to the new HDRI url
#    and flip K2 (Background music Active) on.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "https://raw.githubusercontent.com/Rmoosa2014/vr-tour/main/HDRI01.jpeg"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://raw.githubusercontent.com/Rmoosa2014/vr-tour/main/HDRI01.jpeg") | Out-Null
$ws.Range("K2").Value = $true

# ---------------------------------------------------------------------
# 3. Conditional formatting: duplicate-value highlighting on PanoramaURL
#    now excludes the (changed) F2 cell, and gets its own rule for F24.
# ---------------------------------------------------------------------
$ws.Range("F2:F5").FormatConditions.Delete()
$ws.Range("F3:F5").FormatConditions.AddUniqueValues()
$ws.Range("F3:F5").FormatConditions.Item(1).DupeUnique = 1

$ws.Range("F24").FormatConditions.AddUniqueValues()
$ws.Range("F24").FormatConditions.Item(1).DupeUnique = 1

# ---------------------------------------------------------------------
# 4. Selection / view state, matching the author's final screen position
# ---------------------------------------------------------------------
$ws.Range("M26:R26").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
